# Fixes #687: Parse percentage too when setting datatype to numeric.
# Adds a new "Percentage Text to Number:" example row to the DataTypes
# worksheet, right after the existing "Text to Number:" row, showing a
# text value like "55.12%" converted to the numeric value 0.5512 with a
# percentage number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the old row 34 ("@ format to Number:"), shifting
# everything below it down by one.
$ws.Rows.Item(34).Insert()

# Fill in the new row's label and value.
$ws.Cells.Item(34, 2).Value = "Percentage Text to Number:"
$ws.Cells.Item(34, 3).Value = 0.5512
$ws.Cells.Item(34, 3).NumberFormat = "0.00%"

# Slightly widen column B so the longer label fits (closest value reachable
# through the pixel-quantized ColumnWidth property to the target 25.920625
# raw character width).
$ws.Columns.Item(2).ColumnWidth = 25.17
